{"js": "// Move the \"_GoBack\" bookmark from the end of the \"Submit your self-reflection...\"\n// paragraph into the \"Logic (including ...)\" paragraph (right after the phrase\n// that becomes highlighted), and highlight two phrases in red:\n//   - \"appropriate organization of logic into methods\"\n//   - \"corner cases\"\n\nconst body = context.document.body;\n\n// 1. Remove the bookmark from its current location (it will be re-inserted below).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2. Split \"Logic (including appropriate organization of logic into methods)\"\n//    and highlight the middle phrase red, then re-insert the bookmark right\n//    after it (before the closing parenthesis).\nconst logicHits = body.search(\"appropriate organization of logic into methods\", { matchCase: true });\nlogicHits.load(\"items\");\nawait context.sync();\n\nconst logicPhrase = logicHits.items[0];\nlogicPhrase.font.highlightColor = \"Red\";\nconst afterLogicPhrase = logicPhrase.getRange(\"After\");\nafterLogicPhrase.insertBookmark(\"_GoBack\");\n\n// 3. Split \"Handling all reasonable corner cases\" and highlight \"corner cases\" red.\nconst cornerHits = body.search(\"corner cases\", { matchCase: true });\ncornerHits.load(\"items\");\nawait context.sync();\n\nconst cornerPhrase = cornerHits.items[0];\ncornerPhrase.font.highlightColor = \"Red\";\n\nawait context.sync();\n", "ps1": "# Move the \"_GoBack\" bookmark from the end of the \"Submit your self-reflection...\"\n# paragraph into the \"Logic (including ...)\" paragraph (right after the phrase\n# that becomes highlighted), and highlight two phrases in red:\n#   - \"appropriate organization of logic into methods\"\n#   - \"corner cases\"\n\n$d = $word.ActiveDocument\n\n# 1. Remove the bookmark from its current location (it will be re-inserted below).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Split \"Logic (including appropriate organization of logic into methods)\"\n#    and highlight the middle phrase red, then re-insert the bookmark right\n#    after it (before the closing parenthesis).\n$logicRange = $d.Content\n$logicRange.Find.Execute(\"appropriate organization of logic into methods\") | Out-Null\n$logicRange.Font.HighlightColorIndex = 6\n$bookmarkSpot = $d.Range($logicRange.End, $logicRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkSpot)\n\n# 3. Split \"Handling all reasonable corner cases\" and highlight \"corner cases\" red.\n$cornerRange = $d.Content\n$cornerRange.Find.Execute(\"corner cases\") | Out-Null\n$cornerRange.Font.HighlightColorIndex = 6\n"}
